$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows 46-48 appended to the "Horas" log for the
# "Integración Back-Front - Pruebas con API REST" work.

# Pre-seed the new shared-string values in the same order the original
# author typed the new cell values, so the new unique strings land at
# the same shared-strings indices as in the canonical file.
$ws.Cells.Item(47, 4).Value = "Sprint 3 - Integración BackEnd y FrontEnd"
$ws.Cells.Item(46, 5).Value = "Familiarización son lo creado por Federico"
$ws.Cells.Item(47, 5).Value = "Investigación sobre API REST"

# Row 46 - Bruno Díaz, 23/05/2017, 1 hour, Sprint 3 - BackEnd
$ws.Cells.Item(46, 1).Value = "Bruno Díaz"
$ws.Cells.Item(46, 2).Value = 42878
$ws.Cells.Item(46, 3).Value = 1
$ws.Cells.Item(46, 4).Value = "Sprint 3 - BackEnd"

# Row 47 - Bruno Díaz, 23/05/2017, 1 hour, Sprint 3 - Integración BackEnd y FrontEnd
$ws.Cells.Item(47, 1).Value = "Bruno Díaz"
$ws.Cells.Item(47, 2).Value = 42878
$ws.Cells.Item(47, 3).Value = 1

# Row 48 - Bruno Díaz, 24/05/2017, 2 hours, Sprint 3 - Integración BackEnd y FrontEnd
$ws.Cells.Item(48, 1).Value = "Bruno Díaz"
$ws.Cells.Item(48, 2).Value = 42879
$ws.Cells.Item(48, 3).Value = 2
$ws.Cells.Item(48, 4).Value = "Sprint 3 - Integración BackEnd y FrontEnd"
$ws.Cells.Item(48, 5).Value = "Investigación sobre API REST"

# Apply the same date style used by the rest of the Fecha column (B)
# to the three new date cells, reusing the existing style instead of
# creating a new number-format entry.
$ws.Cells.Item(2, 2).Copy()
$ws.Range("B46:B48").PasteSpecial(-4122)

# Reselect the last touched cell, matching the author's final cursor position
$ws.Range("E39").Select()
